{"js": "// Apply the \"Added many more features\" edits to the Fruitnation review.\n// Each entry is an exact, unique (or uniformly-replaced) original text\n// paired with its replacement. Using body.search with matchCase so we\n// only touch the intended runs, then replacing the whole found range in\n// one shot preserves the run's existing formatting (bold/italic/etc.).\nconst replacements = [\n  [\n    \"Play Fruitnation for Free - Classic Slot Game Review\",\n    \"Play Fruitnation Slot for Free\",\n  ],\n  [\n    \"Straightforward gameplay with no special features\",\n    \"Simple and straightforward gameplay\",\n  ],\n  [\n    \"Excellent RTP of 96.33%\",\n    \"Ability to adjust game settings\",\n  ],\n  [\n    \"Two Gamble modes to choose from\",\n    \"High RTP of 96.33%\",\n  ],\n  [\n    \"Payout level becomes really interesting with the bell and 7 symbols\",\n    \"Option to play with Autostart feature\",\n  ],\n  [\n    \"No free spins or other bonus features\",\n    \"Lack of special features\",\n  ],\n  [\n    \"High maximum bet of \u20ac5,000 may be too steep for some players\",\n    \"High maximum bet of \u20ac5,000\",\n  ],\n  [\n    \"Read our review of Fruitnation, a classic slot game with an excellent RTP of 96.33%. Play for free and explore its straightforward gameplay and Gamble feature.\",\n    \"Read our review of Fruitnation, a classic slot game with straightforward gameplay. Play for free now!\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the Fruitnation review.\n#\n# Strategy: for each (old, new) pair, locate the run containing the exact\n# old text with $d.Content.Find (looping so a pair that occurs more than\n# once - the title - gets every occurrence), then graft a replacement\n# <w:r> onto that exact character range via Range.InsertXML. Unlike\n# Find.Execute(..., Replace:=wdReplaceAll) or Range.Text assignment /\n# Range.Delete (which trigger this host's paragraph run-rebuild and\n# silently drop the neighbouring placeholder \"<w:r/>\" runs the source\n# document uses), grafting InsertXML onto the exact matched range leaves\n# every other run in the paragraph (including empty ones) untouched.\n# The target run's own direct character formatting (e.g. bold/italic) is\n# read back from Range.WordOpenXML first and reapplied verbatim so we\n# don't depend on (possibly style-inherited) Range.Bold/Range.Italic.\n\n$d = $word.ActiveDocument\n\nfunction Escape-XmlText([string]$s) {\n  $s = $s -replace '&', '&amp;'\n  $s = $s -replace '<', '&lt;'\n  $s = $s -replace '>', '&gt;'\n  return $s\n}\n\nfunction Get-DirectRunRPr($rng) {\n  # Pull the literal <w:rPr>...</w:rPr> (if any) straight from the run\n  # that owns $rng's text, via the range's own round-tripped OOXML -\n  # this reflects only direct formatting on the run, not anything the\n  # paragraph/character style contributes.\n  $owx = $rng.WordOpenXML\n  if ($owx -match '(?s)<w:r>(?:<w:rPr>(.*?)</w:rPr>)?<w:t[^>]*>.*?</w:t></w:r>\\s*</w:p>') {\n    return $matches[1]\n  }\n  return $null\n}\n\nfunction Replace-ExactRunText([string]$oldText, [string]$newText) {\n  $replacedCount = 0\n  $searchFrom = 0\n  while ($true) {\n    $scope = $d.Range($searchFrom, $d.Content.End)\n    $find = $scope.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    $found = $find.Execute()\n    if (-not $found) {\n      break\n    }\n\n    $target = $find.Parent\n    $tStart = $target.Start\n    $tEnd = $target.End\n\n    $directRPr = Get-DirectRunRPr $target\n    $rprXml = \"\"\n    if ($directRPr) {\n      $rprXml = \"<w:rPr>\" + $directRPr + \"</w:rPr>\"\n    }\n\n    $escapedNew = Escape-XmlText $newText\n    # Only mark xml:space=\"preserve\" when the replacement text actually has\n    # leading/trailing whitespace that would otherwise be trimmed - matches\n    # how the runs were originally authored (plain <w:t> unless needed).\n    $needsPreserve = ($newText -ne $newText.Trim())\n    $tOpen = \"<w:t>\"\n    if ($needsPreserve) {\n      $tOpen = '<w:t xml:space=\"preserve\">'\n    }\n    $payload = '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n      '<?mso-application progid=\"Word.Document\"?>' +\n      '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n      '<w:body><w:p><w:r>' + $rprXml + $tOpen + $escapedNew + '</w:t></w:r></w:p></w:body>' +\n      '</w:document>' +\n      '</pkg:xmlData></pkg:part></pkg:package>'\n\n    $sub = $d.Range($tStart, $tEnd)\n    $sub.InsertXML($payload)\n\n    $replacedCount = $replacedCount + 1\n    $searchFrom = $tStart + $newText.Length\n  }\n  return $replacedCount\n}\n\nReplace-ExactRunText \"Play Fruitnation for Free - Classic Slot Game Review\" \"Play Fruitnation Slot for Free\" | Out-Null\nReplace-ExactRunText \"Straightforward gameplay with no special features\" \"Simple and straightforward gameplay\" | Out-Null\nReplace-ExactRunText \"Excellent RTP of 96.33%\" \"Ability to adjust game settings\" | Out-Null\nReplace-ExactRunText \"Two Gamble modes to choose from\" \"High RTP of 96.33%\" | Out-Null\nReplace-ExactRunText \"Payout level becomes really interesting with the bell and 7 symbols\" \"Option to play with Autostart feature\" | Out-Null\nReplace-ExactRunText \"No free spins or other bonus features\" \"Lack of special features\" | Out-Null\nReplace-ExactRunText \"High maximum bet of \u20ac5,000 may be too steep for some players\" \"High maximum bet of \u20ac5,000\" | Out-Null\nReplace-ExactRunText \"Read our review of Fruitnation, a classic slot game with an excellent RTP of 96.33%. Play for free and explore its straightforward gameplay and Gamble feature.\" \"Read our review of Fruitnation, a classic slot game with straightforward gameplay. Play for free now!\" | Out-Null\n"}
